$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last_edited_time text for Tháng 8 row (row 7) which shares its
# value with several other rows (D4:D8, D12:D13). All of those rows used the
# same (incorrect) shared-string value, so they all get corrected together.
$ws.Range("D4").Value = "2024-08-24T20:33:00.000Z"
$ws.Range("D5").Value = "2024-08-24T20:33:00.000Z"
$ws.Range("D6").Value = "2024-08-24T20:33:00.000Z"
$ws.Range("D7").Value = "2024-08-24T20:33:00.000Z"
$ws.Range("D8").Value = "2024-08-24T20:33:00.000Z"
$ws.Range("D12").Value = "2024-08-24T20:33:00.000Z"
$ws.Range("D13").Value = "2024-08-24T20:33:00.000Z"

# Update numeric values on row 7 (Tháng 8)
$ws.Range("T7").Value = 27500000
$ws.Range("W7").Value = 94205000
$ws.Range("AA7").Value = 84045000
$ws.Range("AE7").Value = 178250000
$ws.Range("AH7").Value = 129250000
$ws.Range("AK7").Value = 26
$ws.Range("AN7").Value = 49000000
$ws.Range("AQ7").Value = 156750000
